$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 85; this shifts the existing rows
# 85-99 down to 86-100 (and the sheet dimension grows to A1:R100).
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new observation.
$ws.Range("A85").Value = 10
$ws.Range("B85").Value = "Vega Modelo de Temuco"
$ws.Range("C85").Value = "La Araucanía"
$ws.Range("D85").Value = 45173
$ws.Range("E85").Value = 9
$ws.Range("F85").Value = 100112042
$ws.Range("G85").Value = "Locoto"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 80
$ws.Range("K85").Value = 2200
$ws.Range("L85").Value = 2200
$ws.Range("M85").Value = 2200
$ws.Range("N85").Value = "$/kilo"
$ws.Range("O85").Value = "Región de Arica y Parinacota"
$ws.Range("P85").Value = 2200
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = "Hortaliza"
